$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 32.86754772215167
$ws.Range("R2").Value = 295.8079294993651
$ws.Range("S2").Value = 0.006890263497024104
$ws.Range("T2").Value = 0.006890263497024104
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 57.24449264536999
$ws.Range("R3").Value = 515.20043380833
$ws.Range("S3").Value = 0.01200058006804763
$ws.Range("T3").Value = 0.01200058006804763
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 70.71850981348665
$ws.Range("R4").Value = 636.46658832138
$ws.Range("S4").Value = 0.01482523645667073
$ws.Range("T4").Value = 0.01482523645667073
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 734.9671677968547
$ws.Range("R5").Value = 6614.704510171693
$ws.Range("S5").Value = 0.1540765222459479
$ws.Range("T5").Value = 0.1540765222459479
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2683507884157702
$ws.Range("T6").Value = 0.2683507884157702
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 1581.370880109877
$ws.Range("R7").Value = 14232.33792098889
$ws.Range("S7").Value = 0.3315142992287912
$ws.Range("T7").Value = 0.3315142992287912
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 206.9982929729772
$ws.Range("R8").Value = 1862.984636756795
$ws.Range("S8").Value = 0.04339456031448139
$ws.Range("T8").Value = 0.04339456031448138
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 360.5231628433766
$ws.Range("R9").Value = 3244.70846559039
$ws.Range("S9").Value = 0.07557909734461847
$ws.Range("T9").Value = 0.07557909734461847
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 445.3818987876155
$ws.Range("R10").Value = 4008.437089088539
$ws.Range("S10").Value = 0.09336865242864829
$ws.Range("T10").Value = 0.09336865242864829
